# AP110_TestData_PaymentProcessRequest_21C_SubmitPPR.xlsx
# "Add files via upload" - Anu - AP Files Uploaded
#
# Content-level change: on the "Input_Value" sheet, the reference/example
# values that used to live in Q2:S2 (URL / UserName / Password, with Q2
# carrying a hyperlink to the URL) were cleared out, and the now-unused
# hyperlink was removed. The cell formatting (styles) for Q2:S2 is left
# untouched - only the values/hyperlink go away.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")
$ws.Activate()

# Clear the values out of Q2:S2 but keep their existing cell styles.
$ws.Range("Q2:S2").ClearContents()

# Drop the hyperlink that used to live on Q2 (pointed at the URL value).
$ws.Hyperlinks.Delete()

# Reselect so the active selection sits over the cells that changed.
$ws.Range("Q2:S2").Select()
